# Refresh the crypto price/volume table with the latest scraped values.
# (new GitHub Actions run on Tue Aug 27 05:13:56 UTC 2024)
#
# Column D ("Price") cells are switched to Text format ("@") before the
# new value is written so Excel does not auto-coerce numeric-looking
# strings (e.g. "0.370") into numbers, which would silently drop
# significant trailing zeros and change the cell type away from text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell reference -> new display value
$updates = [ordered]@{
    "D2" = "63.048.58"
    "E2" = "  -1.52%  "
    "D3" = "2.684.81"
    "E3" = "  -2.15%  "
    "E4" = "  -0.04%  "
    "D5" = "555.55"
    "E5" = "  -3.03%  "
    "D6" = "158.81"
    "E6" = "  -1.20%  "
    "E7" = "  -0.03%  "
    "D8" = "0.592"
    "E8" = "  -0.74%  "
    "E9" = "  -3.42%  "
    "E10" = "  -2.20%  "
    "D11" = "0.370"
    "E11" = "  -3.88%  "
    "E12" = "  -6.82%  "
    "D13" = "3.159.67"
    "E13" = "  -2.18%  "
    "D14" = "26.47"
    "E14" = "  -1.87%  "
    "D15" = "62.923.56"
    "E15" = "  -1.49%  "
    "E16" = "  -2.17%  "
    "D17" = "2.683.16"
    "E17" = "  -2.33%  "
    "D18" = "11.98"
    "E18" = "  -1.95%  "
    "E19" = "  -3.66%  "
    "D20" = "346.16"
    "E20" = "  -2.41%  "
    "E21" = "  -5.01%  "
    "D22" = "0.999"
    "E22" = "  -0.02%  "
    "D23" = "0.513"
    "E23" = "  -1.87%  "
    "D24" = "63.43"
    "E24" = "  -1.55%  "
    "E25" = "  -1.27%  "
    "E26" = "  +0.11%  "
    "D27" = "8.24"
    "E27" = "  -2.45%  "
    "D28" = "0.0₃0862"
    "E28" = "  -5.84%  "
    "B29" = "Fetch.AI"
    "C29" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D29" = "1.39"
    "E29" = "  +4.26%  "
    "B30" = "Aptos"
    "C30" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D30" = "7.29"
    "E30" = "  +0.75%  "
    "D31" = "1.96"
    "E31" = "  -1.29%  "
    "D32" = "165.26"
    "E32" = "  +0.63%  "
    "D33" = "4.94"
    "E33" = "  +0.42%  "
    "E34" = "  +0.21%  "
    "D36" = "19.53"
    "E36" = "  -3.14%  "
    "E37" = "  -0.85%  "
    "D38" = "350.19"
    "E38" = "  +0.19%  "
    "D39" = "6.35"
    "E39" = "  -0.69%  "
    "D40" = "0.964"
    "E40" = "  -2.61%  "
    "D41" = "4.02"
    "E41" = "  -2.26%  "
    "D42" = "38.48"
    "E42" = "  -0.47%  "
    "D43" = "20.46"
    "E43" = "  -3.39%  "
    "D44" = "20.91"
    "E44" = "  -4.89%  "
    "E45" = "  -1.07%  "
    "E46" = "  -3.94%  "
    "D47" = "0.998"
    "E47" = "  -0.01%  "
    "E48" = "  +0.02%  "
    "D49" = "0.0973"
    "E49" = "  -3.20%  "
    "D50" = "129.62"
    "E50" = "  -3.99%  "
    "D51" = "0.0243"
    "E51" = "  -3.38%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($ref.StartsWith("D")) {
        # Price column: keep it text so formatting (trailing zeros,
        # "." thousands separators, subscript-notation, etc.) survives.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$ref]
}
